$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A, C and D hold text values (residue lists / frame numbers as
# strings) even when they look numeric, so force text format before
# assigning to avoid Excel auto-converting them to numbers.
# (Applied per contiguous column range since multi-area ranges only
# honor the first area for property assignment.)
$ws.Range("A2:A19").NumberFormat = "@"
$ws.Range("C2:C19").NumberFormat = "@"
$ws.Range("D2:D19").NumberFormat = "@"

$ws.Cells.Item(2, 1).Value = "130, 780"
$ws.Cells.Item(2, 2).Value = 1
$ws.Cells.Item(2, 3).Value = "130"
$ws.Cells.Item(2, 4).Value = "3631"

$ws.Cells.Item(3, 1).Value = "455"
$ws.Cells.Item(3, 2).Value = 1
$ws.Cells.Item(3, 3).Value = "455"
$ws.Cells.Item(3, 4).Value = "3666"

$ws.Cells.Item(4, 1).Value = "1073, 1105"
$ws.Cells.Item(4, 2).Value = 1
$ws.Cells.Item(4, 3).Value = "1105"
$ws.Cells.Item(4, 4).Value = "5131"

$ws.Cells.Item(5, 1).Value = "130, 1073"
$ws.Cells.Item(5, 2).Value = 1
$ws.Cells.Item(5, 3).Value = "130"
$ws.Cells.Item(5, 4).Value = "4415"

$ws.Cells.Item(6, 1).Value = "780, 1073, 1105"
$ws.Cells.Item(6, 2).Value = 1
$ws.Cells.Item(6, 3).Value = "1105"
$ws.Cells.Item(6, 4).Value = "4994"

$ws.Cells.Item(7, 1).Value = "130, 455, 780"
$ws.Cells.Item(7, 2).Value = 2
$ws.Cells.Item(7, 3).Value = "130, 130"
$ws.Cells.Item(7, 4).Value = "5269, 6424"

$ws.Cells.Item(8, 1).Value = "130, 1073, 1105"
$ws.Cells.Item(8, 2).Value = 1
$ws.Cells.Item(8, 3).Value = "1105"
$ws.Cells.Item(8, 4).Value = "5399"

$ws.Cells.Item(9, 1).Value = "423, 748, 780, 1073"
$ws.Cells.Item(9, 2).Value = 2
$ws.Cells.Item(9, 3).Value = "780, 780"
$ws.Cells.Item(9, 4).Value = "5677, 5887"

$ws.Cells.Item(10, 1).Value = "423, 1105, 1105"
$ws.Cells.Item(10, 2).Value = 1
$ws.Cells.Item(10, 3).Value = "1105"
$ws.Cells.Item(10, 4).Value = "5331"

$ws.Cells.Item(11, 1).Value = "423, 1073, 1105"
$ws.Cells.Item(11, 2).Value = 1
$ws.Cells.Item(11, 3).Value = "1105"
$ws.Cells.Item(11, 4).Value = "5433"

$ws.Cells.Item(12, 1).Value = "98, 130, 748, 780, 1073"
$ws.Cells.Item(12, 2).Value = 1
$ws.Cells.Item(12, 3).Value = "130"
$ws.Cells.Item(12, 4).Value = "6016"

$ws.Cells.Item(13, 1).Value = "98, 130, 748, 1073"
$ws.Cells.Item(13, 2).Value = 2
$ws.Cells.Item(13, 3).Value = "130, 130"
$ws.Cells.Item(13, 4).Value = "5582, 6488"

$ws.Cells.Item(14, 1).Value = "423, 748, 780, 1073, SF"
$ws.Cells.Item(14, 2).Value = 1
$ws.Cells.Item(14, 3).Value = "780"
$ws.Cells.Item(14, 4).Value = "6202"

$ws.Cells.Item(15, 1).Value = "130, 423, 748, 1073"
$ws.Cells.Item(15, 2).Value = 1
$ws.Cells.Item(15, 3).Value = "130"
$ws.Cells.Item(15, 4).Value = "6561"

$ws.Cells.Item(16, 1).Value = "130, 780, 1073, 1105"
$ws.Cells.Item(16, 2).Value = 1
$ws.Cells.Item(16, 3).Value = "1105"
$ws.Cells.Item(16, 4).Value = "6359"

$ws.Cells.Item(17, 1).Value = "98, 130, 130, 423, 780"
$ws.Cells.Item(17, 2).Value = 1
$ws.Cells.Item(17, 3).Value = "130"
$ws.Cells.Item(17, 4).Value = "6727"

$ws.Cells.Item(18, 1).Value = "98, 130, 423, 1073"
$ws.Cells.Item(18, 2).Value = 1
$ws.Cells.Item(18, 3).Value = "130"
$ws.Cells.Item(18, 4).Value = "6670"

$ws.Cells.Item(19, 1).Value = "98, 130, 455, 780"
$ws.Cells.Item(19, 2).Value = 1
$ws.Cells.Item(19, 3).Value = "780"
$ws.Cells.Item(19, 4).Value = "6748"
